# Auto-generated Excel COM-interop script to update cryptocurrency price/volume data
# per commit "Updated symbol list on Sun Jan  8 08:45:27 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target cells hold plain text (price / percentage strings), not numbers.
# Force text number-format before writing so COM does not coerce the
# strings into floating point numbers (price) or percentage-typed
# numbers (volume), which would lose exact formatting/precision.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "260.98"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.14%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.62%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.714"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.07%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06223"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.31%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.730"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.75%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8493"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.42%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9119"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.24%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1402"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.15%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04930"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.88%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07079"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.39%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03082"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.40%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.20%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001540"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.37%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006155"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.25%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005958"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.92%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.446"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.19%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.173"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.76%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.168"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.04%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1311"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.39%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.117"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.62%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04248"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.21%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001196"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.91%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004070"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.06%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.05%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.08%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03930"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.55%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.24%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004134"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.01%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.39%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01330"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-18.67%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.84%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2484"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "83.59%"
